# Adicionando EPL 17_18 18_19
# Adds three new team/abbreviation rows (Stoke/STO, West Brom/WBA, Hull/HUL)
# to the de/para table, plus two trailing "shell" rows that carry formatting
# only (no values) - mirrors the author's commit that appended the 2017/18
# and 2018/19 Premier League relegated-club mappings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows ---------------------------------------------------
$ws.Range("A53").Value = "Stoke"
$ws.Range("B53").Value = "STO"

$ws.Range("A54").Value = "West Brom"
$ws.Range("B54").Value = "WBA"

$ws.Range("A55").Value = "Hull"
$ws.Range("B55").Value = "HUL"
$ws.Rows.Item(55).RowHeight = 15

# --- Trailing formatted-but-empty rows (matches the source workbook) -
$ws.Range("A56").HorizontalAlignment = -4108
$ws.Range("A56").WrapText = $true

$ws.Range("B57").WrapText = $true

# --- View state: scroll down, zoom in, select the new blank row ------
$excel.ActiveWindow.Zoom = 110
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("A56:B56").Select()
